# Adds season-record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF ---
# Copy the existing header formatting (bold, centered, bordered) from A1
# onto the new header cells before setting their values.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2 through 51): season record for each player's team ---
# Wins = 80, Losses = 82, Ties = 0
$lastRow = 51
$ws.Range("AD2:AD$lastRow").Value = 80
$ws.Range("AE2:AE$lastRow").Value = 82
$ws.Range("AF2:AF$lastRow").Value = 0

Write-Host "Season record columns (Wins/Losses/Ties) added for rows 1-$lastRow"
